$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.316.70"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "3.687.00"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "680.68"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").Value = "159.26"
$ws.Range("E6").Value = "  -2.21%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("E10").Value = "  -3.74%  "
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("D12").Value = "0.0000232"
$ws.Range("E12").Value = "  -3.25%  "
$ws.Range("D13").Value = "4.307.49"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").Value = "32.49"
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("D15").Value = "3.694.59"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "69.302.94"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("E17").Value = "  +1.97%  "
$ws.Range("D18").Value = "16.08"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").Value = "6.45"
$ws.Range("E19").Value = "  -2.42%  "
$ws.Range("D20").Value = "468.34"
$ws.Range("E20").Value = "  -2.61%  "
$ws.Range("D21").Value = "9.95"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("D23").Value = "79.83"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").Value = "3.833.39"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  -5.79%  "
$ws.Range("E27").Value = "  -4.13%  "
$ws.Range("D28").Value = "9.13"
$ws.Range("E28").Value = "  -4.00%  "
$ws.Range("D29").Value = "2.69"
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("E30").Value = "  -4.44%  "
$ws.Range("D31").Value = "6.64"
$ws.Range("E31").Value = "  -2.90%  "
$ws.Range("E32").Value = "  -3.36%  "
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("D34").Value = "26.93"
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("D35").Value = "3.675.57"
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("E36").Value = "  -5.81%  "
$ws.Range("D37").Value = "8.29"
$ws.Range("E37").Value = "  -2.39%  "
$ws.Range("D38").Value = "6.26"
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("E40").Value = "  -2.01%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("E42").Value = "  -2.75%  "
$ws.Range("D43").Value = "170.28"
$ws.Range("E43").Value = "  +4.02%  "
$ws.Range("D44").Value = "0.942"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("D45").Value = "47.59"
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("D46").Value = "28.53"
$ws.Range("E46").Value = "  -5.37%  "
$ws.Range("E47").Value = "  -2.67%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.18%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.38%  "
$ws.Range("D50").Value = "0.000276"
$ws.Range("E50").Value = "  -3.65%  "
$ws.Range("D51").Value = "7.79"
$ws.Range("E51").Value = "  -3.79%  "
